$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.076.22'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.306.48'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +3.96%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.67'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').Value = '2.665.25'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '2.299.72'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.788'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '42.984.99'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.59%  '
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E30').Value = '  -7.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.15%  '
$ws.Range('E35').Value = '  +5.56%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('E40').Value = '  +2.52%  '
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '2.010.25'
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').Value = '2.530.55'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.32%  '
$ws.Range('E51').Value = '  +1.72%  '
